$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Support:" row -> placeholder token departmentFull becomes "support"
#    The original run layout is: "${" | <spellStart/> "departmentFull" <spellEnd/> | "}"
#    The target layout keeps the "${" and "}" runs as-is, drops the spell-check
#    proofErr wrappers, and swaps the middle run's text for "support".
# ---------------------------------------------------------------------------
$label = $d.Content
$null = $label.Find.Execute('Support:')

# Replace the whole placeholder text; Word merges the old "${" / "departmentFull"
# / "}" runs (and drops their now-stale proofErr spell-check markers) into one
# new run.
$scope = $d.Range($label.End, $d.Content.End)
$null = $scope.Find.Execute('${departmentFull}', $true, $false, $false, $false, $false, `
                             $true, 1, $false, '${support}', 2)

# Re-split "support" back out into its own run (distinct from "${" and "}") by
# forcing a direct-formatting transition on just that word, then clearing it
# back so no stray formatting override is left behind.
$merged = $d.Content
$null = $merged.Find.Execute('${support}')
$word1 = $d.Range($merged.Start + 2, $merged.End - 1)
$word1.Bold = $true
$word1.Bold = $false

# ---------------------------------------------------------------------------
# 2)-4) These placeholders are already alone in their cell/paragraph, so a
#    straight Find & Replace of the full "${...}" token collapses the three
#    runs ("${" / name / "}") plus their proofErr wrappers into a single run.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute('${ilos}', $true, $false, $false, $false, $false, `
                                 $true, 1, $false, '${ilos}', 2)

$null = $d.Content.Find.Execute('${budgetSource}', $true, $false, $false, $false, $false, `
                                 $true, 1, $false, '${budgetSource}', 2)

$null = $d.Content.Find.Execute('${sig_sscp}', $true, $false, $false, $false, $false, `
                                 $true, 1, $false, '${sig_sscp}', 2)

$null = $d.Content.Find.Execute('${sig_dean}', $true, $false, $false, $false, $false, `
                                 $true, 1, $false, '${sig_dean}', 2)
